$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows for the new "August Red" / "August pearl" weekly price records
# (existing rows 536:569 shift down to 542:575)
$ws.Rows("536:541").Insert()

# Row 536: August Red / Especial
$ws.Range("A536").Value = 5
$ws.Range("B536").Value = 'Macroferia Regional de Talca'
$ws.Range("C536").Value = 'Maule'
$ws.Range("D536").Value = 44610
$ws.Range("E536").Value = 7
$ws.Range("F536").Value = 'Fruta'
$ws.Range("G536").Value = 100103
$ws.Range("H536").Value = 'Frutos de hueso (carozo)'
$ws.Range("I536").Value = 100103006
$ws.Range("J536").Value = 'Nectarín'
$ws.Range("K536").Value = 'August Red'
$ws.Range("L536").Value = 'Especial'
$ws.Range("M536").Value = 300
$ws.Range("N536").Value = 12000
$ws.Range("O536").Value = 12000
$ws.Range("P536").Value = 12000
$ws.Range("Q536").Value = '$/bandeja 18 kilos granel'
$ws.Range("R536").Value = 'Región de O''Higgins'
$ws.Range("S536").Value = 667
$ws.Range("T536").Value = 18

# Row 537: August Red / Extra (doble especial)
$ws.Range("A537").Value = 5
$ws.Range("B537").Value = 'Macroferia Regional de Talca'
$ws.Range("C537").Value = 'Maule'
$ws.Range("D537").Value = 44610
$ws.Range("E537").Value = 7
$ws.Range("F537").Value = 'Fruta'
$ws.Range("G537").Value = 100103
$ws.Range("H537").Value = 'Frutos de hueso (carozo)'
$ws.Range("I537").Value = 100103006
$ws.Range("J537").Value = 'Nectarín'
$ws.Range("K537").Value = 'August Red'
$ws.Range("L537").Value = 'Extra (doble especial)'
$ws.Range("M537").Value = 300
$ws.Range("N537").Value = 14000
$ws.Range("O537").Value = 14000
$ws.Range("P537").Value = 14000
$ws.Range("Q537").Value = '$/bandeja 18 kilos granel'
$ws.Range("R537").Value = 'Región de O''Higgins'
$ws.Range("S537").Value = 778
$ws.Range("T537").Value = 18

# Row 538: August Red / Primera
$ws.Range("A538").Value = 5
$ws.Range("B538").Value = 'Macroferia Regional de Talca'
$ws.Range("C538").Value = 'Maule'
$ws.Range("D538").Value = 44610
$ws.Range("E538").Value = 7
$ws.Range("F538").Value = 'Fruta'
$ws.Range("G538").Value = 100103
$ws.Range("H538").Value = 'Frutos de hueso (carozo)'
$ws.Range("I538").Value = 100103006
$ws.Range("J538").Value = 'Nectarín'
$ws.Range("K538").Value = 'August Red'
$ws.Range("L538").Value = 'Primera'
$ws.Range("M538").Value = 200
$ws.Range("N538").Value = 10000
$ws.Range("O538").Value = 10000
$ws.Range("P538").Value = 10000
$ws.Range("Q538").Value = '$/bandeja 18 kilos granel'
$ws.Range("R538").Value = 'Región de O''Higgins'
$ws.Range("S538").Value = 556
$ws.Range("T538").Value = 18

# Row 539: August pearl / Especial
$ws.Range("A539").Value = 5
$ws.Range("B539").Value = 'Macroferia Regional de Talca'
$ws.Range("C539").Value = 'Maule'
$ws.Range("D539").Value = 44610
$ws.Range("E539").Value = 7
$ws.Range("F539").Value = 'Fruta'
$ws.Range("G539").Value = 100103
$ws.Range("H539").Value = 'Frutos de hueso (carozo)'
$ws.Range("I539").Value = 100103006
$ws.Range("J539").Value = 'Nectarín'
$ws.Range("K539").Value = 'August pearl'
$ws.Range("L539").Value = 'Especial'
$ws.Range("M539").Value = 300
$ws.Range("N539").Value = 12000
$ws.Range("O539").Value = 12000
$ws.Range("P539").Value = 12000
$ws.Range("Q539").Value = '$/bandeja 18 kilos granel'
$ws.Range("R539").Value = 'Región de O''Higgins'
$ws.Range("S539").Value = 667
$ws.Range("T539").Value = 18

# Row 540: August pearl / Extra (doble especial)
$ws.Range("A540").Value = 5
$ws.Range("B540").Value = 'Macroferia Regional de Talca'
$ws.Range("C540").Value = 'Maule'
$ws.Range("D540").Value = 44610
$ws.Range("E540").Value = 7
$ws.Range("F540").Value = 'Fruta'
$ws.Range("G540").Value = 100103
$ws.Range("H540").Value = 'Frutos de hueso (carozo)'
$ws.Range("I540").Value = 100103006
$ws.Range("J540").Value = 'Nectarín'
$ws.Range("K540").Value = 'August pearl'
$ws.Range("L540").Value = 'Extra (doble especial)'
$ws.Range("M540").Value = 300
$ws.Range("N540").Value = 14000
$ws.Range("O540").Value = 14000
$ws.Range("P540").Value = 14000
$ws.Range("Q540").Value = '$/bandeja 18 kilos granel'
$ws.Range("R540").Value = 'Región de O''Higgins'
$ws.Range("S540").Value = 778
$ws.Range("T540").Value = 18

# Row 541: August pearl / Primera
$ws.Range("A541").Value = 5
$ws.Range("B541").Value = 'Macroferia Regional de Talca'
$ws.Range("C541").Value = 'Maule'
$ws.Range("D541").Value = 44610
$ws.Range("E541").Value = 7
$ws.Range("F541").Value = 'Fruta'
$ws.Range("G541").Value = 100103
$ws.Range("H541").Value = 'Frutos de hueso (carozo)'
$ws.Range("I541").Value = 100103006
$ws.Range("J541").Value = 'Nectarín'
$ws.Range("K541").Value = 'August pearl'
$ws.Range("L541").Value = 'Primera'
$ws.Range("M541").Value = 200
$ws.Range("N541").Value = 10000
$ws.Range("O541").Value = 10000
$ws.Range("P541").Value = 10000
$ws.Range("Q541").Value = '$/bandeja 18 kilos granel'
$ws.Range("R541").Value = 'Región de O''Higgins'
$ws.Range("S541").Value = 556
$ws.Range("T541").Value = 18
